# Add 20 more tugboats (NB023-NB040 / CP0023-CP0040) to the schedule.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tugboat rows (24-41): TugBoatId, CurrentStatus, CaptainId, StartWorkingTime, EndWorkingTime
$newRows = @(
    @("NB023", "Maintenance", "CP0023"),
    @("NB024", "Free",        "CP0024"),
    @("NB025", "Free",        "CP0025"),
    @("NB026", "Free",        "CP0026"),
    @("NB027", "Free",        "CP0027"),
    @("NB028", "Free",        "CP0028"),
    @("NB029", "Free",        "CP0029"),
    @("NB030", "Maintenance", "CP0030"),
    @("NB031", "Free",        "CP0031"),
    @("NB032", "Free",        "CP0032"),
    @("NB033", "Free",        "CP0033"),
    @("NB034", "Free",        "CP0034"),
    @("NB035", "Free",        "CP0035"),
    @("NB036", "Free",        "CP0036"),
    @("NB037", "Free",        "CP0037"),
    @("NB038", "Free",        "CP0038"),
    @("NB039", "Free",        "CP0039"),
    @("NB040", "Free",        "CP0040")
)

$startTime = 0.33333333333333298
$endTime = 0.75

$row = 24
foreach ($entry in $newRows) {
    # Copy the formatting of the last existing data row (23) down onto the
    # new row so the date-formatted D/E cells keep the same style (s="5").
    $ws.Range("A23:E23").Copy()
    $ws.Range("A" + $row + ":E" + $row).PasteSpecial(-4122)

    $ws.Range("A" + $row).Value2 = $entry[0]
    $ws.Range("B" + $row).Value2 = $entry[1]
    $ws.Range("C" + $row).Value2 = $entry[2]
    $ws.Range("D" + $row).Value2 = $startTime
    $ws.Range("E" + $row).Value2 = $endTime

    $row = $row + 1
}

# Trailing partially-formatted blank rows (42-44): only B (no style) and
# D/E (date style carried over) were touched, A and C stay untouched.
$ws.Range("B41").Copy()
$ws.Range("B42:B44").PasteSpecial(-4122)
$ws.Range("D41:E41").Copy()
$ws.Range("D42:E44").PasteSpecial(-4122)

# Restore the view: scroll so row 13 is at the top and the active cell is G27.
$ws.Range("G27").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
